$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) VALOR MORA total (E11) and Cant. Periodos (F13) go up because a
#    new period (2509) is being added to the account statement.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 650700
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 2) The "Novedad de Ingreso" / "Novedad de Retiro" header columns
#    (H15/I15) were swapped.
# ------------------------------------------------------------------
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# ------------------------------------------------------------------
# 3) Insert three new detail rows (21:23) for period 2509, right
#    after the existing 2508 rows and before the signature block.
# ------------------------------------------------------------------
$formatSourceNormal = $ws.Range("B19:J19")
$formatSourceBottom = $ws.Range("B20:J20")

$ws.Rows("21:23").Insert()

# Row 20 is no longer the last row of the table, give it the normal
# interior-row formatting; the brand-new row 23 becomes the new last
# row and gets the special bottom-border formatting.
$formatSourceNormal.Copy()
$ws.Range("B20:J22").PasteSpecial(-4122)
$formatSourceBottom.Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Fill in the values for the new 2509 period rows, following the
#    same worker ordering already used for 2507/2508.
# ------------------------------------------------------------------
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1001835668"
$ws.Range("D21").Value = "JORGE ANDRES GARCIA OSORIO"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1017182514"
$ws.Range("D22").Value = "ISABEL CRISTINA ALZATE MARIN"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "92642113"
$ws.Range("D23").Value = "RICK PETER HERNANDEZ RUSSO"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 122000
$ws.Range("G23").Value = 3050000
